$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 77027.30499999999
$ws.Range("I11").Value = 77027.30499999999
$ws.Range("K11").Value = 77027.30499999999
$ws.Range("M11").Value = -76887.30499999999

$ws.Range("H12").Value = 38806.08
$ws.Range("I12").Value = 248
$ws.Range("K12").Value = 248
$ws.Range("M12").Value = -78

$ws.Range("H33").Value = 1245.1875
$ws.Range("J33").Value = 2687
$ws.Range("L33").Value = 2687
$ws.Range("N33").Value = -3145

$ws.Range("H40").Value = 1681.5769
$ws.Range("I40").Value = 1615.1
$ws.Range("J40").Value = 1903.1666
$ws.Range("K40").Value = 1615.1
$ws.Range("L40").Value = 1903.1666
$ws.Range("M40").Value = -1440.1
$ws.Range("N40").Value = -2253.1666

$ws.Range("H64").Value = 252697.5
$ws.Range("I64").Value = 1000000
$ws.Range("J64").Value = 3596.6667
$ws.Range("K64").Value = 1000000
$ws.Range("L64").Value = 3596.6667
$ws.Range("M64").Value = -999752
$ws.Range("N64").Value = -4092.6667

$ws.Range("H67").Value = 252697.5
$ws.Range("I67").Value = 1000000
$ws.Range("J67").Value = 3596.6667
$ws.Range("K67").Value = 1000000
$ws.Range("L67").Value = 3596.6667
$ws.Range("M67").Value = -999142
$ws.Range("N67").Value = -5312.6667

$ws.Range("H74").Value = 6499.857
$ws.Range("I74").Value = 5266.6665
$ws.Range("J74").Value = 7424.75
$ws.Range("K74").Value = 5266.6665
$ws.Range("L74").Value = 7424.75
$ws.Range("M74").Value = -4330.6665
$ws.Range("N74").Value = -9296.75

$ws.Range("H77").Value = 6499.857
$ws.Range("I77").Value = 5266.6665
$ws.Range("J77").Value = 7424.75
$ws.Range("K77").Value = 26333.3325
$ws.Range("L77").Value = 37123.75
$ws.Range("M77").Value = -21653.3325
$ws.Range("N77").Value = -46483.75

$ws.Range("H116").Value = 1695.5
$ws.Range("I116").Value = 1219.4445
$ws.Range("K116").Value = 1219.4445
$ws.Range("M116").Value = 2222.5555

$ws.Range("H129").Value = 923.2632
$ws.Range("J129").Value = 965.4167
$ws.Range("L129").Value = 2896.2501
$ws.Range("N129").Value = -12896.2501

$ws.Range("H132").Value = 10876112
$ws.Range("I132").Value = 11911742
$ws.Range("J132").Value = 1994.5
$ws.Range("K132").Value = 35735226
$ws.Range("L132").Value = 5983.5
$ws.Range("M132").Value = -35732696
$ws.Range("N132").Value = -11043.5

$ws.Range("H137").Value = 1308.9272
$ws.Range("I137").Value = 828.7
$ws.Range("J137").Value = 1583.3429
$ws.Range("K137").Value = 2486.1
$ws.Range("L137").Value = 4750.028700000001
$ws.Range("M137").Value = 63.89999999999964
$ws.Range("N137").Value = -9850.028700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1649.1333
$ws.Range("I45").Value = 1734.2858
$ws.Range("K45").Value = 1734.2858
$ws.Range("M45").Value = -1357.2858

$ws.Range("H132").Value = 4786.909
$ws.Range("I132").Value = 4850.778
$ws.Range("K132").Value = 14552.334
$ws.Range("M132").Value = -12022.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 155731.69
$ws.Range("I105").Value = 126686.25
$ws.Range("J105").Value = 202204.4
$ws.Range("K105").Value = 126686.25
$ws.Range("L105").Value = 202204.4
$ws.Range("M105").Value = -124939.25
$ws.Range("N105").Value = -205698.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 532.8889
$ws.Range("I22").Value = 338.8
$ws.Range("J22").Value = 775.5
$ws.Range("K22").Value = 338.8
$ws.Range("L22").Value = 775.5
$ws.Range("M22").Value = 11.19999999999999
$ws.Range("N22").Value = -1475.5

$ws.Range("H31").Value = 17428.27
$ws.Range("I31").Value = 28770.027
$ws.Range("J31").Value = 2305.926
$ws.Range("K31").Value = 28770.027
$ws.Range("L31").Value = 2305.926
$ws.Range("M31").Value = -28475.027
$ws.Range("N31").Value = -2895.926

$ws.Range("H34").Value = 17428.27
$ws.Range("I34").Value = 28770.027
$ws.Range("J34").Value = 2305.926
$ws.Range("K34").Value = 28770.027
$ws.Range("L34").Value = 2305.926
$ws.Range("M34").Value = -28568.027
$ws.Range("N34").Value = -2709.926

$ws.Range("H50").Value = 10594.286
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 10594.286
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 10594.286
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -11844.286

$ws.Range("H59").Value = 24805
$ws.Range("J59").Value = 24805
$ws.Range("L59").Value = 24805
$ws.Range("N59").Value = -27095

$ws.Range("H60").Value = 18073.334
$ws.Range("I60").Value = 12000
$ws.Range("J60").Value = 19288
$ws.Range("K60").Value = 12000
$ws.Range("L60").Value = 19288
$ws.Range("M60").Value = -11489
$ws.Range("N60").Value = -20310

$ws.Range("H86").Value = 2552.4
$ws.Range("I86").Value = 1990.8334
$ws.Range("J86").Value = 3070.7693
$ws.Range("K86").Value = 1990.8334
$ws.Range("L86").Value = 3070.7693
$ws.Range("M86").Value = -867.8334
$ws.Range("N86").Value = -5316.7693

$ws.Range("H89").Value = 2552.4
$ws.Range("I89").Value = 1990.8334
$ws.Range("J89").Value = 3070.7693
$ws.Range("K89").Value = 9954.166999999999
$ws.Range("L89").Value = 15353.8465
$ws.Range("M89").Value = -4338.166999999999
$ws.Range("N89").Value = -26585.8465

$ws.Range("H105").Value = 998.8570999999999
$ws.Range("J105").Value = 1199.75
$ws.Range("L105").Value = 1199.75
$ws.Range("N105").Value = -4693.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 721.9722
$ws.Range("I5").Value = 561.86365
$ws.Range("J5").Value = 973.5714
$ws.Range("K5").Value = 1685.59095
$ws.Range("L5").Value = 2920.7142
$ws.Range("M5").Value = -1573.59095
$ws.Range("N5").Value = -3144.7142

$ws.Range("H55").Value = 7283.8096
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 7283.8096
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21851.4288
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -22205.4288

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws.Range("H131").Value = 2302.3914
$ws.Range("I131").Value = 730
$ws.Range("J131").Value = 2373.8635
$ws.Range("K131").Value = 2190
$ws.Range("L131").Value = 7121.5905
$ws.Range("M131").Value = 2850
$ws.Range("N131").Value = -17201.5905

$ws.Range("H132").Value = 2635.4736
$ws.Range("I132").Value = 1849.8334
$ws.Range("J132").Value = 2998.077
$ws.Range("K132").Value = 16648.5006
$ws.Range("L132").Value = 26982.693
$ws.Range("M132").Value = -14118.5006
$ws.Range("N132").Value = -32042.693

$ws.Range("H135").Value = 721.9722
$ws.Range("I135").Value = 561.86365
$ws.Range("J135").Value = 973.5714
$ws.Range("K135").Value = 5056.77285
$ws.Range("L135").Value = 8762.142600000001
$ws.Range("M135").Value = -2521.77285
$ws.Range("N135").Value = -13832.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 8949.5
$ws.Range("I31").Value = 415.5
$ws.Range("J31").Value = 26017.5
$ws.Range("K31").Value = 415.5
$ws.Range("L31").Value = 26017.5
$ws.Range("M31").Value = -123.5
$ws.Range("N31").Value = -26601.5

$ws.Range("H37").Value = 8949.5
$ws.Range("I37").Value = 415.5
$ws.Range("J37").Value = 26017.5
$ws.Range("K37").Value = 415.5
$ws.Range("L37").Value = 26017.5
$ws.Range("M37").Value = -138.5
$ws.Range("N37").Value = -26571.5

$ws.Range("H68").Value = 35995
$ws.Range("J68").Value = 35995
$ws.Range("L68").Value = 35995
$ws.Range("N68").Value = -37617

$ws.Range("H71").Value = 35995
$ws.Range("J71").Value = 35995
$ws.Range("L71").Value = 107985
$ws.Range("N71").Value = -116097

$ws.Range("H136").Value = 24099.555
$ws.Range("J136").Value = 24099.555
$ws.Range("L136").Value = 72298.66500000001
$ws.Range("N136").Value = -77398.66500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3333
$ws.Range("I46").Value = 998.5
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 998.5
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -810.5
$ws.Range("N46").Value = -4376

$ws.Range("H94").Value = 24900
$ws.Range("J94").Value = 24900
$ws.Range("L94").Value = 24900
$ws.Range("N94").Value = -26252

$ws.Range("H119").Value = 34042.5
$ws.Range("J119").Value = 34042.5
$ws.Range("L119").Value = 34042.5
$ws.Range("N119").Value = -43718.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 9366.5
$ws.Range("I14").Value = 849.75
$ws.Range("K14").Value = 849.75
$ws.Range("M14").Value = -681.75

$ws.Range("H100").Value = 71761.36
$ws.Range("I100").Value = 125193.625
$ws.Range("J100").Value = 518.3333
$ws.Range("K100").Value = 250387.25
$ws.Range("L100").Value = 1036.6666
$ws.Range("M100").Value = -249846.25
$ws.Range("N100").Value = -2118.6666
